$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.997.20'
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').Value = '2.086.70'
$ws.Range('E3').Value = '  +2.66%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''228.80'
$ws.Range('E5').Value = '  +0.68%  '
$ws.Range('E6').Value = '  +1.15%  '
$ws.Range('D7').Value = '''60.72'
$ws.Range('E7').Value = '  +0.90%  '
$ws.Range('D9').Value = '''0.381'
$ws.Range('E9').Value = '  -0.26%  '
$ws.Range('D10').Value = '''0.0838'
$ws.Range('E10').Value = '  +2.47%  '
$ws.Range('E11').Value = '  -0.12%  '
$ws.Range('D12').Value = '2.396.73'
$ws.Range('E12').Value = '  +2.56%  '
$ws.Range('D13').Value = '''14.61'
$ws.Range('E13').Value = '  +0.41%  '
$ws.Range('E14').Value = '  +3.29%  '
$ws.Range('D15').Value = '''5.50'
$ws.Range('E15').Value = '  +6.59%  '
$ws.Range('D16').Value = '''0.772'
$ws.Range('E16').Value = '  +1.52%  '
$ws.Range('D17').Value = '2.105.63'
$ws.Range('E17').Value = '  +3.66%  '
$ws.Range('D18').Value = '37.562.28'
$ws.Range('E18').Value = '  -0.73%  '
$ws.Range('E19').Value = '  +2.27%  '
$ws.Range('D20').Value = '''70.03'
$ws.Range('E20').Value = '  +0.22%  '
$ws.Range('D21').Value = '0.0₃0837'
$ws.Range('E21').Value = '  +1.31%  '
$ws.Range('D22').Value = '''224.03'
$ws.Range('E22').Value = '  -0.22%  '
$ws.Range('E23').Value = '  +0.62%  '
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('E25').Value = '  +2.74%  '
$ws.Range('D26').Value = '''169.54'
$ws.Range('E26').Value = '  +1.49%  '
$ws.Range('D27').Value = '''9.38'
$ws.Range('E27').Value = '  +0.78%  '
$ws.Range('D28').Value = '''0.133'
$ws.Range('E28').Value = '  +3.69%  '
$ws.Range('D29').Value = '''18.96'
$ws.Range('E29').Value = '  +0.31%  '
$ws.Range('E30').Value = '  +3.79%  '
$ws.Range('E31').Value = '  -0.35%  '
$ws.Range('D32').Value = '''2.37'
$ws.Range('E32').Value = '  +10.72%  '
$ws.Range('D33').Value = '''4.43'
$ws.Range('E33').Value = '  +0.61%  '
$ws.Range('D34').Value = '''4.65'
$ws.Range('D35').Value = '''0.0606'
$ws.Range('E35').Value = '  +0.08%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').Value = '''2.41'
$ws.Range('E36').Value = '  +6.19%  '
$ws.Range('B37').Value = 'THORChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D37').Value = '''6.47'
$ws.Range('E37').Value = '  +0.55%  '
$ws.Range('E38').Value = '  +8.43%  '
$ws.Range('D39').Value = '''1.00'
$ws.Range('E39').Value = '  -0.09%  '
$ws.Range('D40').Value = '''17.95'
$ws.Range('E40').Value = '  +4.21%  '
$ws.Range('D41').Value = '1.546.69'
$ws.Range('E41').Value = '  +1.48%  '
$ws.Range('D42').Value = '''100.08'
$ws.Range('E42').Value = '  +4.20%  '
$ws.Range('E43').Value = '  +0.44%  '
$ws.Range('E44').Value = '  -0.10%  '
$ws.Range('D45').Value = '''0.0904'
$ws.Range('E45').Value = '  -1.20%  '
$ws.Range('D46').Value = '''4.15'
$ws.Range('E46').Value = '  +3.80%  '
$ws.Range('E47').Value = '  +0.80%  '
$ws.Range('E48').Value = '  +1.11%  '
$ws.Range('D49').Value = '''2.99'
$ws.Range('E49').Value = '  +1.00%  '
$ws.Range('D50').Value = '''7.21'
$ws.Range('E50').Value = '  +1.46%  '
$ws.Range('D51').Value = '2.283.46'
$ws.Range('E51').Value = '  +2.58%  '
